# Update "江西-漫展信息.xlsx" (Jiangxi con-info workbook) to the data
# snapshot generated at commit 456a3b4 for the gh-pages output.
#
# Two sheets carry data changes:
#   "展览"    (Exhibitions) - rows 1-32, gains a new row 33
#   "全部类型" (All types)   - rows 1-33, gains a new row 34 (it has one
#                              extra pre-existing row vs "展览" because it
#                              also includes the single "演出" / Performance
#                              entry at its row 5)
#
# For both sheets: a handful of "想去人数" (interest-count, column F)
# values tick upward, and a brand-new con ("南昌·第一届异次元动漫嘉年华")
# is appended as the last row.

$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet "展览" (Exhibitions): rows 1-32 -> +row 33
# ======================================================================
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(3, 6).Value = 186
$ws1.Cells.Item(4, 6).Value = 182
$ws1.Cells.Item(5, 6).Value = 5111
$ws1.Cells.Item(7, 6).Value = 48
$ws1.Cells.Item(8, 6).Value = 18
$ws1.Cells.Item(9, 6).Value = 565
$ws1.Cells.Item(10, 6).Value = 522
$ws1.Cells.Item(12, 6).Value = 19
$ws1.Cells.Item(13, 6).Value = 1425
$ws1.Cells.Item(14, 6).Value = 3817
$ws1.Cells.Item(15, 6).Value = 421
$ws1.Cells.Item(16, 6).Value = 152
$ws1.Cells.Item(17, 6).Value = 136
$ws1.Cells.Item(18, 6).Value = 89
$ws1.Cells.Item(19, 6).Value = 3064
$ws1.Cells.Item(20, 6).Value = 143
$ws1.Cells.Item(21, 6).Value = 623
$ws1.Cells.Item(25, 6).Value = 76
$ws1.Cells.Item(28, 6).Value = 65
$ws1.Cells.Item(29, 6).Value = 283
$ws1.Cells.Item(31, 6).Value = 49
$ws1.Cells.Item(32, 6).Value = 5

# Append row 33, copying column A's bold/bordered/centered style from
# the row above it (row 32).
$ws1.Range("A32").Copy() | Out-Null
$ws1.Range("A33").PasteSpecial(-4122) | Out-Null

$ws1.Cells.Item(33, 1).Value = 32
# Leading apostrophe forces text (matches the other "开始时间" cells,
# which are stored as literal strings, not date serials).
$ws1.Cells.Item(33, 2).Value = "'2024-08-06"
$ws1.Cells.Item(33, 3).Value = "南昌·第一届异次元动漫嘉年华"
$ws1.Cells.Item(33, 4).Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws1.Cells.Item(33, 5).Value = "2024.08.06 08:00-08.06 17:00"
$ws1.Cells.Item(33, 6).Value = 0
$ws1.Cells.Item(33, 7).Value = 40
$ws1.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws1.Cells.Item(33, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/0W8gIOYx1712764727806.jpeg"

# ======================================================================
# Sheet "全部类型" (All types): rows 1-33 -> +row 34
# ======================================================================
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(3, 6).Value = 186
$ws4.Cells.Item(4, 6).Value = 182
$ws4.Cells.Item(6, 6).Value = 5111
$ws4.Cells.Item(8, 6).Value = 48
$ws4.Cells.Item(9, 6).Value = 18
$ws4.Cells.Item(10, 6).Value = 565
$ws4.Cells.Item(11, 6).Value = 522
$ws4.Cells.Item(13, 6).Value = 19
$ws4.Cells.Item(14, 6).Value = 1425
$ws4.Cells.Item(15, 6).Value = 3817
$ws4.Cells.Item(16, 6).Value = 421
$ws4.Cells.Item(17, 6).Value = 152
$ws4.Cells.Item(18, 6).Value = 136
$ws4.Cells.Item(19, 6).Value = 89
$ws4.Cells.Item(20, 6).Value = 3064
$ws4.Cells.Item(21, 6).Value = 143
$ws4.Cells.Item(22, 6).Value = 623
$ws4.Cells.Item(26, 6).Value = 76
$ws4.Cells.Item(29, 6).Value = 65
$ws4.Cells.Item(30, 6).Value = 283
$ws4.Cells.Item(32, 6).Value = 49
$ws4.Cells.Item(33, 6).Value = 5

# Append row 34, copying column A's bold/bordered/centered style from
# the row above it (row 33).
$ws4.Range("A33").Copy() | Out-Null
$ws4.Range("A34").PasteSpecial(-4122) | Out-Null

$ws4.Cells.Item(34, 1).Value = 33
# Leading apostrophe forces text (matches the other "开始时间" cells,
# which are stored as literal strings, not date serials).
$ws4.Cells.Item(34, 2).Value = "'2024-08-06"
$ws4.Cells.Item(34, 3).Value = "南昌·第一届异次元动漫嘉年华"
$ws4.Cells.Item(34, 4).Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws4.Cells.Item(34, 5).Value = "2024.08.06 08:00-08.06 17:00"
$ws4.Cells.Item(34, 6).Value = 0
$ws4.Cells.Item(34, 7).Value = 40
$ws4.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws4.Cells.Item(34, 9).Value = "//i0.hdslb.com/bfs/openplatform/202404/0W8gIOYx1712764727806.jpeg"

Write-Output "applied con-info update (456a3b4)"
